$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Bump the cached "datetimeFigureOut" date placeholder text from
#    8/5/2020 -> 8/6/2020 everywhere it is cached: the slide master and
#    every slide layout's Date Placeholder shape.
# ---------------------------------------------------------------------------
$masterShapes = $p.SlideMaster.Shapes
for ($i = 1; $i -le $masterShapes.Count; $i++) {
    $shp = $masterShapes.Item($i)
    if ($shp.Type -eq 14) {
        if ($shp.PlaceholderFormat.Type -eq 16) {
            $shp.TextFrame.TextRange.Text = "8/6/2020"
        }
    }
}

$layouts = $p.SlideMaster.CustomLayouts
for ($k = 1; $k -le $layouts.Count; $k++) {
    $lytShapes = $layouts.Item($k).Shapes
    for ($i = 1; $i -le $lytShapes.Count; $i++) {
        $shp = $lytShapes.Item($i)
        if ($shp.Type -eq 14) {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $shp.TextFrame.TextRange.Text = "8/6/2020"
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 6 ("Support during the Workshop"): credit Franny Buderman
#    alongside Brent Pease and John Yeiser.
# ---------------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
$contentShape = $slide6.Shapes.Item(2)
$tr = $contentShape.TextFrame.TextRange

$full = $tr.Text
$idx = $full.IndexOf("Brent Pease and John ")
$target = $tr.Characters($idx + 1, 21)
$target.Text = "Brent Pease, Franny Buderman, and John "

$full2 = $tr.Text
$idxB = $full2.IndexOf("Buderman")
$nameRun = $tr.Characters($idxB + 1, 8)
$nameRun.Text = "Buderman"
